# Auto-generated: update crypto price/volume table values to match the
# "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.732.55'
$ws.Range("E2").Value = '  -6.31%  '
$ws.Range("D3").Value = '3.295.96'
$ws.Range("E3").Value = '  -7.45%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.48%  '
$ws.Range("D5").Value = "'182.21"
$ws.Range("E5").Value = '  -11.09%  '
$ws.Range("D6").Value = "'519.07"
$ws.Range("E6").Value = '  -7.11%  '
$ws.Range("D7").Value = "'0.598"
$ws.Range("E7").Value = '  -1.05%  '
$ws.Range("D8").Value = '3.296.06'
$ws.Range("E8").Value = '  -7.27%  '
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("D10").Value = "'0.622"
$ws.Range("E10").Value = '  -7.06%  '
$ws.Range("D11").Value = "'59.06"
$ws.Range("E11").Value = '  -5.52%  '
$ws.Range("D12").Value = "'0.132"
$ws.Range("E12").Value = '  -8.82%  '
$ws.Range("D13").Value = "'0.0000256"
$ws.Range("E13").Value = '  -7.57%  '
$ws.Range("D14").Value = "'9.15"
$ws.Range("E14").Value = '  -8.46%  '
$ws.Range("D15").Value = '3.798.18'
$ws.Range("E15").Value = '  -8.35%  '
$ws.Range("E16").Value = '  -5.35%  '
$ws.Range("D17").Value = '3.280.91'
$ws.Range("E17").Value = '  -8.37%  '
$ws.Range("D18").Value = "'17.73"
$ws.Range("E18").Value = '  -6.25%  '
$ws.Range("D19").Value = '63.432.56'
$ws.Range("E19").Value = '  -6.59%  '
$ws.Range("D20").Value = "'11.01"
$ws.Range("E20").Value = '  -8.56%  '
$ws.Range("D21").Value = "'0.950"
$ws.Range("E21").Value = '  -9.60%  '
$ws.Range("D22").Value = "'373.13"
$ws.Range("E22").Value = '  -5.17%  '
$ws.Range("D23").Value = "'11.33"
$ws.Range("E23").Value = '  -6.45%  '
$ws.Range("D24").Value = "'80.37"
$ws.Range("E24").Value = '  -3.98%  '
$ws.Range("D25").Value = "'3.68"
$ws.Range("E25").Value = '  -10.06%  '
$ws.Range("D26").Value = "'3.86"
$ws.Range("E26").Value = '  +1.36%  '
$ws.Range("E27").Value = '  -2.48%  '
$ws.Range("D28").Value = "'2.65"
$ws.Range("E28").Value = '  -6.49%  '
$ws.Range("D29").Value = "'11.43"
$ws.Range("E29").Value = '  -7.03%  '
$ws.Range("D30").Value = "'8.33"
$ws.Range("E30").Value = '  -7.25%  '
$ws.Range("D31").Value = "'28.62"
$ws.Range("E31").Value = '  -7.95%  '
$ws.Range("D32").Value = "'648.80"
$ws.Range("E32").Value = '  -9.75%  '
$ws.Range("D33").Value = "'6.71"
$ws.Range("E33").Value = '  -9.82%  '
$ws.Range("D34").Value = "'11.20"
$ws.Range("E34").Value = '  -5.99%  '
$ws.Range("D35").Value = "'59.42"
$ws.Range("E35").Value = '  -6.59%  '
$ws.Range("D36").Value = "'0.105"
$ws.Range("E36").Value = '  -5.19%  '
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("D38").Value = "'0.389"
$ws.Range("E38").Value = '  -6.30%  '
$ws.Range("D39").Value = "'36.25"
$ws.Range("E39").Value = '  -10.96%  '
$ws.Range("D40").Value = "'0.995"
$ws.Range("E40").Value = '  -0.48%  '
$ws.Range("D41").Value = '2.997.05'
$ws.Range("E41").Value = '  -4.60%  '
$ws.Range("D42").Value = "'0.125"
$ws.Range("E42").Value = '  -4.08%  '
$ws.Range("D43").Value = '0.0₃0652'
$ws.Range("E43").Value = '  -9.61%  '
$ws.Range("D44").Value = "'2.68"
$ws.Range("E44").Value = '  -16.22%  '
$ws.Range("D45").Value = "'2.43"
$ws.Range("E45").Value = '  -4.56%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = "'0.0389"
$ws.Range("E46").Value = '  -4.19%  '
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").Value = "'2.59"
$ws.Range("E47").Value = '  -5.74%  '
$ws.Range("D48").Value = "'2.79"
$ws.Range("E48").Value = '  +4.85%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = "'0.126"
$ws.Range("E49").Value = '  -2.67%  '
$ws.Range("B50").Value = 'ApeXProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D50").Value = "'2.94"
$ws.Range("E50").Value = '  -3.64%  '
$ws.Range("D51").Value = "'2.47"
$ws.Range("E51").Value = '  -20.28%  '
